# Update "paises.xlsx" (sheet "Pais") with a refreshed data pull.
#
# The country list is sorted descending by "Casos totales" (col B). Some
# countries' updated totals changed their rank, so their row swapped with
# a neighboring row (names below are re-labelled in place to reflect the
# new sort order); other rows just received refreshed statistics while
# keeping their country. Finally the "last updated" timestamp is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose rank (and therefore country label) changed -----------------
$ws.Range("A43").Value  = "Emiratos Arabes Unidos"
$ws.Range("B43").Value  = 110039
$ws.Range("C43").Value  = 1431
$ws.Range("D43").Value  = 101659
$ws.Range("E43").Value  = 7930
$ws.Range("G43").Value  = 2
$ws.Range("H43").Value  = 450

$ws.Range("A44").Value  = "Kazajistan"
$ws.Range("B44").Value  = 108984
$ws.Range("C44").Value  = 83
$ws.Range("D44").Value  = 104346
$ws.Range("E44").Value  = 2870
$ws.Range("H44").Value  = 1768

$ws.Range("A78").Value  = "Tunez"
$ws.Range("B78").Value  = 34790
$ws.Range("C78").Value  = 2234
$ws.Range("D78").Value  = 5032
$ws.Range("E78").Value  = 29246
$ws.Range("G78").Value  = 34
$ws.Range("H78").Value  = 512

$ws.Range("A79").Value  = "Dinamarca"
$ws.Range("B79").Value  = 33101
$ws.Range("D79").Value  = 27225
$ws.Range("E79").Value  = 5202
$ws.Range("H79").Value  = 674

$ws.Range("A88").Value  = "Eslovaquia"
$ws.Range("B88").Value  = 22296
$ws.Range("C88").Value  = 1410
$ws.Range("D88").Value  = 6709
$ws.Range("E88").Value  = 15521
$ws.Range("G88").Value  = 5
$ws.Range("H88").Value  = 66

$ws.Range("A89").Value  = "Croacia"
$ws.Range("B89").Value  = 21741
$ws.Range("C89").Value  = 748
$ws.Range("D89").Value  = 18197
$ws.Range("E89").Value  = 3210
$ws.Range("G89").Value  = 4
$ws.Range("H89").Value  = 334

$ws.Range("A90").Value  = "Camerun"
$ws.Range("B90").Value  = 21203
$ws.Range("D90").Value  = 20117
$ws.Range("E90").Value  = 663
$ws.Range("H90").Value  = 423

$ws.Range("A91").Value  = "Republica de Macedonia"
$ws.Range("B91").Value  = 21193
$ws.Range("D91").Value  = 16397
$ws.Range("E91").Value  = 3996
$ws.Range("H91").Value  = 800

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

$ws.Range("A217").Value = "Montserrat"
$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1

# --- Rows that only got refreshed statistics (country unchanged) -----------
$ws.Range("B4").Value   = 8093600
$ws.Range("C4").Value   = 3347
$ws.Range("D4").Value   = 5227279
$ws.Range("E4").Value   = 2645421
$ws.Range("G4").Value   = 27
$ws.Range("H4").Value   = 220900

$ws.Range("B19").Value  = 382959
$ws.Range("C19").Value  = 1684
$ws.Range("D19").Value  = 297449
$ws.Range("E19").Value  = 79917
$ws.Range("G19").Value  = 16
$ws.Range("H19").Value  = 5593

$ws.Range("B32").Value  = 164477
$ws.Range("C32").Value  = 4016
$ws.Range("D32").Value  = 122714
$ws.Range("E32").Value  = 36162
$ws.Range("G32").Value  = 66
$ws.Range("H32").Value  = 5601

$ws.Range("B45").Value  = 107776
$ws.Range("C45").Value  = 563
$ws.Range("D45").Value  = 93908
$ws.Range("E45").Value  = 12807
$ws.Range("G45").Value  = 8
$ws.Range("H45").Value  = 1061

$ws.Range("B58").Value  = 68704
$ws.Range("C58").Value  = 2823
$ws.Range("E58").Value  = 16798
$ws.Range("G58").Value  = 1
$ws.Range("H58").Value  = 2106

$ws.Range("D127").Value = 4932
$ws.Range("E127").Value = 165

$ws.Range("B151").Value = 2942
$ws.Range("C151").Value = 102
$ws.Range("E151").Value = 1576

$ws.Range("B175").Value = 575
$ws.Range("C175").Value = 10
$ws.Range("E175").Value = 31

# --- Footer timestamp --------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 12:29"
